# Generate Report for handoff
# This script updates the localization-status workbook:
#  - rotates the "current handoff" file from one generated UUID.md to a new one
#  - adds a "Handoff transform failed" row for a second generated UUID.md
#  - keeps the ".localization-config" row, now pushed down a row
# applied uniformly across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldUuid  = "2f0da9c7-7a5c-4630-b477-8813e138fd33"
$newUuid  = "5ac20b70-63ad-459b-9428-d49f65bd972c"
$newUuid2 = "088b3ce5-cb0f-4af8-8006-c0a0c5144da3"
$oldHash  = "087236e3e78bc56066c7383767f17546e21850de"
$newHash  = "cf4ccd8494ceb9bf3bb47a99591c2f5d31db45e6"

$newMdUrl      = "https://github.com/OpenLocalizationTest/oltest/blob/5908d66dee621f7e6534eb8f3c11b5ee1d85c166/e2e/$newUuid.md"
$newMd2Url     = "https://github.com/OpenLocalizationTest/oltest/blob/5908d66dee621f7e6534eb8f3c11b5ee1d85c166/e2e/$newUuid2.md"
$localConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/5908d66dee621f7e6534eb8f3c11b5ee1d85c166/.localization-config"
$newZhCnXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e4b11d177293f3dd6e9fea56f82520f223ea5fca/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$newUuid.$newHash.zh-cn.xlf"
$newDeDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ab4a5dc3e67c57d29d4b97f0a342c4e68dbec218/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$newUuid.$newHash.de-de.xlf"

$newZhCnTimestamp = "2016-02-16 15:19:34"
$newDeDeTimestamp = "2016-02-16 15:19:51"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Preserve the hyperlink-style formatting (font color + underline) on the
# rows that already have it, then also stamp it onto the new row before
# clearing + rebuilding the hyperlinks collection.
$ws1.Range("A3").Copy()
$ws1.Range("A4").PasteSpecial()
$ws1.Range("B3").Copy()
$ws1.Range("B4").PasteSpecial()
$ws1.Range("C3").Copy()
$ws1.Range("C4").PasteSpecial()

# New row 4 = old row 3 content (.localization-config / Not to be localized)
$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

# Row 2 now points at the new uuid markdown file
$ws1.Range("A2").Value = $newUuid + ".md"
# Row 3 becomes the "Handoff transform failed" entry for the 2nd uuid
$ws1.Range("A3").Value = $newUuid2 + ".md"
$ws1.Range("B3").Value = "Handoff transform failed"
$ws1.Range("C3").Value = "Handoff transform failed"

# Hyperlinks.Delete() on any range nukes the whole sheet's collection in
# this engine, so rebuild it from scratch in final order.
$ws1.Range("A1").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $newMdUrl, "", "", $newUuid + ".md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $newMd2Url, "", "", $newUuid2 + ".md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), $localConfigUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Stamp the new row 4 with the same styles as row 3 before overwriting values
$ws2.Range("A3").Copy()
$ws2.Range("A4").PasteSpecial()
$ws2.Range("B3").Copy()
$ws2.Range("B4").PasteSpecial()
$ws2.Range("D3").Copy()
$ws2.Range("D4").PasteSpecial()
$ws2.Range("G3").Copy()
$ws2.Range("G4").PasteSpecial()
$ws2.Range("H3").Copy()
$ws2.Range("H4").PasteSpecial()

# New row 4 = old row 3 content
$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

# Row 2: rotate to the new uuid / xlf / timestamp
$ws2.Range("A2").Value = $newUuid + ".md"
$ws2.Range("C2").Value = $newUuid + "." + $newHash + ".zh-cn.xlf"
$ws2.Range("D2").Value = $newZhCnTimestamp

# Row 3: becomes "Handoff transform failed" for the 2nd uuid
$ws2.Range("A3").Value = $newUuid2 + ".md"
$ws2.Range("B3").Value = "Handoff transform failed"
$ws2.Range("D3").Value = "0001-01-01 00:00:00"
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Ignored"

# Re-assert the number format on every datetime cell so the shared style
# stays valid for all of them (old + new) after the stylesheet is rewritten.
$ws2.Range("D2").NumberFormat = $dateFmt
$ws2.Range("D3").NumberFormat = $dateFmt
$ws2.Range("D4").NumberFormat = $dateFmt

$ws2.Range("A1").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $newMdUrl, "", "", $newUuid + ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $newZhCnXlfUrl, "", "", $newUuid + "." + $newHash + ".zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $newMd2Url, "", "", $newUuid2 + ".md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), $localConfigUrl, "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A3").Copy()
$ws3.Range("A4").PasteSpecial()
$ws3.Range("B3").Copy()
$ws3.Range("B4").PasteSpecial()
$ws3.Range("D3").Copy()
$ws3.Range("D4").PasteSpecial()
$ws3.Range("G3").Copy()
$ws3.Range("G4").PasteSpecial()
$ws3.Range("H3").Copy()
$ws3.Range("H4").PasteSpecial()

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Range("A2").Value = $newUuid + ".md"
$ws3.Range("C2").Value = $newUuid + "." + $newHash + ".de-de.xlf"
$ws3.Range("D2").Value = $newDeDeTimestamp

$ws3.Range("A3").Value = $newUuid2 + ".md"
$ws3.Range("B3").Value = "Handoff transform failed"
$ws3.Range("D3").Value = "0001-01-01 00:00:00"
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Ignored"

$ws3.Range("D2").NumberFormat = $dateFmt
$ws3.Range("D3").NumberFormat = $dateFmt
$ws3.Range("D4").NumberFormat = $dateFmt

$ws3.Range("A1").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $newMdUrl, "", "", $newUuid + ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $newDeDeXlfUrl, "", "", $newUuid + "." + $newHash + ".de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $newMd2Url, "", "", $newUuid2 + ".md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), $localConfigUrl, "", "", ".localization-config") | Out-Null
